$d = $word.ActiveDocument

# Locate the "9.3.5" run (part of "IntelliJ IDEA 2019.3.5, JUnit 5.4.2").
$anchor = $d.Content
$anchor.Find.Execute("9.3.5", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$anchorStart = $anchor.Start
$anchorEnd = $anchor.End

# Temporarily perturb "9.3.5"'s formatting so the upcoming insertion is not
# silently coalesced into it (adjacent runs sharing identical direct
# formatting get merged by the engine's paragraph-rebuild pass).
$rng935 = $d.Range($anchorStart, $anchorEnd)
$rng935.Font.Size = 9

# Insert the new "IU" build-version text right after "9.3.5".
$newText = " IU 183.5912.21"
$ins = $d.Range($anchorEnd, $anchorEnd)
$ins.InsertAfter($newText)
$newStart = $anchorEnd
$newEnd = $anchorEnd + $newText.Length

# Restore "9.3.5"'s original formatting (Bold, navy color, size 10pt) using a
# precise, non-collapsed, single-run range so it doesn't trigger another
# merge pass.
$fix935 = $d.Range($anchorStart, $anchorEnd)
$fix935.Font.Bold = $true
$fix935.Font.Color = 8388608
$fix935.Font.Size = 10

# Apply the same formatting to the newly inserted text.
$fixNew = $d.Range($newStart, $newEnd)
$fixNew.Font.Bold = $true
$fixNew.Font.Color = 8388608
$fixNew.Font.Size = 10
